$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header value (X-1P-ZUUL-HOST=...) from F2:F4, leaving the cells empty
$ws.Range("F2:F4").ClearContents()
